# Slide 39 ("Reference via NPM" -> "See the Code"): retitle and add
# the "how to use it" steps to the body placeholder.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(39)

# --- Title shape: "Reference via NPM" -> "See the Code" ---
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "See the Code"

# --- Content placeholder: add the three bullet paragraphs ---
$bodyRange = $s.Shapes.Item(2).TextFrame.TextRange
# The placeholder starts as a single, empty paragraph (just a paragraph
# end mark). Insert the new paragraphs ahead of that mark so it is kept
# as the trailing paragraph end of the final ("Provide "/"the service")
# bullet, instead of being replaced outright.
[void]$bodyRange.InsertBefore("Install the module" + [char]13 + "Import the module/service" + [char]13 + "Provide ")
# Split the final bullet into two runs ("Provide " + "the service") by
# appending the second run right before that preserved paragraph mark.
[void]$bodyRange.InsertAfter("the service")
